# ------------------------------------------------------------------
# country_comparison figs: 'support' -> 'accept' wording update,
# plus refreshed survey rows/values (per commit diff).
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A labels (rows 2-9) ---
$ws.Range("A2").Value = "Accepts tax on world top 1% to finance global poverty reduction`n(Additional 15% tax on income over [`$120k/year in PPP])"
$ws.Range("A3").Value = "Accepts tax on world top 3% to finance global poverty reduction`n(Additional 15% tax over [`$80k], 30% over [`$120k], 45% over [`$1M])"
$ws.Range("A4").Value = "Prefers sustainable future"
$ws.Range("A5").Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""
$ws.Range("A6").Value = "Would support a global movement to tackle CC, tax millionaires,`n and fund LICs (either petition, demonstrate, strike, or donate)"
$ws.Range("A7").Value = "More likely to vote for party if part of worldwide`ncoalition for climate action and global redistribution"
$ws.Range("A8").Value = "Accepts reparations for colonization and slavery in`nthe form of funding education and technology transfers"
$ws.Range("A9").Value = "`"My taxes should go towards solving global problems`""

# --- Numeric grid B2:N9 ---
$grid = New-Object 'object[,]' 8,13
$grid[0,0] = 0.692894784662911
$grid[0,1] = 0.709953011533533
$grid[0,2] = 0.694020071682012
$grid[0,3] = 0.716874917652575
$grid[0,4] = 0.817931401759025
$grid[0,5] = 0.688331328542526
$grid[0,6] = 0.727696489623426
$grid[0,7] = 0.680287179370093
$grid[0,8] = 0.611054424765204
$grid[0,9] = 0.667587247975338
$grid[0,10] = 0.734266953673365
$grid[0,11] = 0.816793423425975
$grid[0,12] = 0.613078913946868
$grid[1,0] = 0.648431801604668
$grid[1,1] = 0.654418197725284
$grid[1,2] = 0.69867423579194
$grid[1,3] = 0.631062611744521
$grid[1,4] = 0.71347266445345
$grid[1,5] = 0.698563809102872
$grid[1,6] = 0.668248725994209
$grid[1,7] = 0.688695322348303
$grid[1,8] = 0.412823621757237
$grid[1,9] = 0.548789112211414
$grid[1,10] = 0.745575769472929
$grid[1,11] = 0.820236659115239
$grid[1,12] = 0.587193769163202
$grid[2,0] = 0.664778686805119
$grid[2,1] = 0.680960854092527
$grid[2,2] = 0.701612903225806
$grid[2,3] = 0.683610867659947
$grid[2,4] = 0.728285077951002
$grid[2,5] = 0.568822553897181
$grid[2,6] = 0.726299694189602
$grid[2,7] = 0.667752442996743
$grid[2,8] = 0.65607476635514
$grid[2,9] = 0.706269349845201
$grid[2,10] = 0.685934489402698
$grid[2,11] = 0.668763102725367
$grid[2,12] = 0.602535832414553
$grid[3,0] = 0.705793226381462
$grid[3,1] = 0.7602300376023
$grid[3,2] = 0.762917933130699
$grid[3,3] = 0.740149094781683
$grid[3,4] = 0.856145251396648
$grid[3,5] = 0.825301204819277
$grid[3,6] = 0.830601092896175
$grid[3,7] = 0.652413793103448
$grid[3,8] = 0.658256880733945
$grid[3,9] = 0.663065843621399
$grid[3,10] = 0.765895953757225
$grid[3,11] = 0.917480998914224
$grid[3,12] = 0.555436337625179
$grid[4,0] = 0.611777124330845
$grid[4,1] = 0.672953736654804
$grid[4,2] = 0.668202764976959
$grid[4,3] = 0.670464504820333
$grid[4,4] = 0.736080178173719
$grid[4,5] = 0.648424543946932
$grid[4,6] = 0.723241590214067
$grid[4,7] = 0.642779587404995
$grid[4,8] = 0.598130841121495
$grid[4,9] = 0.506191950464396
$grid[4,10] = 0.506191950464396
$grid[4,11] = 0.560447239692523
$grid[4,12] = 0.612458654906284
$grid[5,0] = 0.664188137644821
$grid[5,1] = 0.711714770797963
$grid[5,2] = 0.710578842315369
$grid[5,3] = 0.705192629815745
$grid[5,4] = 0.801272507913065
$grid[5,5] = 0.646090534979424
$grid[5,6] = 0.769662921348315
$grid[5,7] = 0.693877551020408
$grid[5,8] = 0.58695652173913
$grid[5,9] = 0.515331355093966
$grid[5,10] = $null
$grid[5,11] = $null
$grid[5,12] = 0.669950738916256
$grid[6,0] = 0.461188014718766
$grid[6,1] = 0.502795031055901
$grid[6,2] = 0.436893203883495
$grid[6,3] = 0.442563482466747
$grid[6,4] = 0.685958024097665
$grid[6,5] = $null
$grid[6,6] = 0.511201629327902
$grid[6,7] = 0.461068702290076
$grid[6,8] = $null
$grid[6,9] = $null
$grid[6,10] = $null
$grid[6,11] = $null
$grid[6,12] = 0.407318053880177
$grid[7,0] = 0.622476446837147
$grid[7,1] = 0.612353567625133
$grid[7,2] = 0.441176470588235
$grid[7,3] = 0.62962962962963
$grid[7,4] = 0.771253333873262
$grid[7,5] = 0.641833810888252
$grid[7,6] = 0.712018140589569
$grid[7,7] = 0.584084084084084
$grid[7,8] = 0.532786885245902
$grid[7,9] = 0.601431980906921
$grid[7,10] = 0.577994428969359
$grid[7,11] = 0.882267286664075
$grid[7,12] = 0.571005917159763
$ws.Range("B2:N9").Value = $grid
